$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the blank spacer row (old row 3, which only held an empty B3 cell).
#    This shifts: old row4 ("Area" label) -> row3, old row5 (year headers) -> row4,
#    old row6 (area values) -> row5.
$ws.Rows.Item(3).Delete()

# 2. Drop the "(according to the population census data)" subtitle text from A2,
#    leaving row 2 as a blank spacer row.
$ws.Cells.Item(2, 1).ClearContents()

# 3. Clear the now stray wrap-text placeholder cell in B1 (row 1 should only have A1).
$ws.Cells.Item(1, 2).ClearContents()

# 4. Keep only the last data column (2014 / 799.53) - remove the 1989 and 2002
#    columns (B and C), so the former column D becomes column B.
$ws.Columns.Item(3).Delete()
$ws.Columns.Item(2).Delete()

# 5. The remaining two data rows (years header + values) now get an explicit,
#    larger row height.
$ws.Rows.Item(4).RowHeight = 20.1
$ws.Rows.Item(5).RowHeight = 20.1
